# Swap the two theme colour palettes used by this deck.
#
# Before: the Slide Master's theme ("theme2.xml") used the "Integral"
# palette, while the Notes Master's theme ("theme1.xml") used the plain
# "Office Theme" palette.
#
# After:  the Slide Master's theme now carries the "Office Theme" palette
# and the Notes Master's theme carries the "Integral" palette - i.e. the
# two palettes have effectively traded places.
#
# We drive this through the DrawingML 12-colour theme scheme exposed on a
# slide (it resolves to the deck's single Slide Master theme), writing the
# RGB value of every slot (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette: the plain "Office Theme" colours (previously used only
# by the Notes Master) now become the Slide Master's theme colours.
$tcs.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
